# Apply the diary (Дневник) updates: add three new date/activity rows
# (16 марта, 17 марта, 18 марта) into rows 21-23, and update the current
# selection/scroll position to reflect where the user ended up editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 21: 16  марта (note: two spaces between "16" and "марта" in source)
$ws.Range("A21").Value = "16  марта"
$ws.Range("B21").Value = "Написание потока для мониторинга работы стадий конвейера"

# Row 22: 17 марта
$ws.Range("A22").Value = "17 марта"
$ws.Range("B22").Value = "Написание теорической части вкр, оптимизация передачи данных между стадиями конвейера (оптимизация копирования)"

# Row 23: 18 марта
$ws.Range("A23").Value = "18 марта"
$ws.Range("B23").Value = "Усовершенствование потока для мониторинга конвейера: добавлено снятие замеров времени, которое стадии затрачивают на выполнение функтора"

# Match column B's alignment to column A's existing style (xlRight) so the
# new cells carry the same cell style index as the rest of the sheet.
$ws.Range("B21:B23").HorizontalAlignment = -4152

# Reflect the final selection / scrolled viewport from the saved workbook.
$ws.Range("B24").Select()
$excel.ActiveWindow.ScrollRow = 7
